$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 88.5
$ws.Range("I9").Value = 32.75
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 32.75
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 136.25
$ws.Range("N9").Value = -538

$ws.Range("H121").Value = 3110.625
$ws.Range("J121").Value = 3982.5
$ws.Range("L121").Value = 11947.5
$ws.Range("N121").Value = -15441.5

$ws.Range("H137").Value = 2115.5293
$ws.Range("I137").Value = 1905.8
$ws.Range("J137").Value = 2415.1428
$ws.Range("K137").Value = 5717.4
$ws.Range("L137").Value = 7245.428400000001
$ws.Range("M137").Value = -3167.4
$ws.Range("N137").Value = -12345.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1920.4667
$ws.Range("I2").Value = 1800.1
$ws.Range("J2").Value = 2161.2
$ws.Range("K2").Value = 1800.1
$ws.Range("L2").Value = 2161.2
$ws.Range("M2").Value = -1687.1
$ws.Range("N2").Value = -2387.2

$ws.Range("H116").Value = 1920.4667
$ws.Range("I116").Value = 1800.1
$ws.Range("J116").Value = 2161.2
$ws.Range("K116").Value = 1800.1
$ws.Range("L116").Value = 2161.2
$ws.Range("M116").Value = 493.9000000000001
$ws.Range("N116").Value = -6749.2

$ws.Range("H128").Value = 45975
$ws.Range("J128").Value = 45975
$ws.Range("L128").Value = 45975
$ws.Range("N128").Value = -55935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1920.4667
$ws.Range("I3").Value = 1800.1
$ws.Range("J3").Value = 2161.2
$ws.Range("K3").Value = 1800.1
$ws.Range("L3").Value = 2161.2
$ws.Range("M3").Value = -1686.1
$ws.Range("N3").Value = -2389.2

$ws.Range("H97").Value = 7777
$ws.Range("I97").Value = 7777
$ws.Range("K97").Value = 7777
$ws.Range("M97").Value = -6786

$ws.Range("H105").Value = 3750
$ws.Range("I105").Value = 3750
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3750
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = -2003
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 32295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 32295
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = 32295
$ws.Range("N68").Value = -33793
$ws.Range("L68").ClearContents()

$ws.Range("H71").Value = 32295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 32295
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = 96885
$ws.Range("N71").Value = -104373
$ws.Range("L71").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 70
$ws.Range("I6").Value = 70
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 210
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = -97
$ws.Range("M6").ClearContents()

$ws.Range("H13").Value = 1889.8334
$ws.Range("I13").Value = 1101.3334
$ws.Range("J13").Value = 2678.3333
$ws.Range("K13").Value = 3304.0002
$ws.Range("L13").Value = 8034.999899999999
$ws.Range("M13").Value = -3136.0002
$ws.Range("N13").Value = -8370.999899999999

$ws.Range("H34").Value = 1318.5714
$ws.Range("J34").Value = 1641.8182
$ws.Range("L34").Value = 4925.4546
$ws.Range("N34").Value = -5093.4546

$ws.Range("H36").Value = 2120.4
$ws.Range("I36").Value = 150.5
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 451.5
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = -282.5
$ws.Range("N36").Value = -30338

$ws.Range("H46").Value = 459.2

$ws.Range("H48").Value = 1896.6666
$ws.Range("J48").Value = 1896.6666
$ws.Range("L48").Value = 5689.9998
$ws.Range("N48").Value = -6189.9998

$ws.Range("H55").Value = 7502.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 7502.5
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = 22507.5
$ws.Range("N55").Value = -22861.5
$ws.Range("L55").ClearContents()

$ws.Range("H59").Value = 1010
$ws.Range("I59").Value = 630
$ws.Range("J59").Value = 1580
$ws.Range("K59").Value = 1890
$ws.Range("L59").Value = 4740
$ws.Range("M59").Value = -1350
$ws.Range("N59").Value = -5820

$ws.Range("H69").Value = 475.63635
$ws.Range("J69").Value = 459
$ws.Range("L69").Value = 1377
$ws.Range("N69").Value = -2999

$ws.Range("H72").Value = 475.63635
$ws.Range("J72").Value = 459
$ws.Range("L72").Value = 4131
$ws.Range("N72").Value = -12243

$ws.Range("H75").Value = 877.5
$ws.Range("I75").Value = 483.33334
$ws.Range("J75").Value = 1271.6666
$ws.Range("K75").Value = 1450.00002
$ws.Range("L75").Value = 3814.9998
$ws.Range("M75").Value = -452.0000199999999
$ws.Range("N75").Value = -5810.9998

$ws.Range("H78").Value = 877.5
$ws.Range("I78").Value = 483.33334
$ws.Range("J78").Value = 1271.6666
$ws.Range("K78").Value = 4350.00006
$ws.Range("L78").Value = 11444.9994
$ws.Range("M78").Value = 641.9999399999997
$ws.Range("N78").Value = -21428.9994

$ws.Range("H81").Value = 4776
$ws.Range("I81").Value = 313
$ws.Range("J81").Value = 7007.5
$ws.Range("K81").Value = 939
$ws.Range("L81").Value = 21022.5
$ws.Range("M81").Value = 184
$ws.Range("N81").Value = -23268.5

$ws.Range("H84").Value = 4776
$ws.Range("I84").Value = 313
$ws.Range("J84").Value = 7007.5
$ws.Range("K84").Value = 2817
$ws.Range("L84").Value = 63067.5
$ws.Range("M84").Value = 2799
$ws.Range("N84").Value = -74299.5

$ws.Range("H94").Value = 6800
$ws.Range("J94").Value = 6800
$ws.Range("L94").Value = 20400
$ws.Range("N94").Value = -21752

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H112").Value = 1623.6364
$ws.Range("J112").Value = 1623.6364
$ws.Range("L112").Value = 4870.9092
$ws.Range("N112").Value = -7086.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2526.2856
$ws.Range("I16").Value = 1590
$ws.Range("J16").Value = 2900.8
$ws.Range("K16").Value = 1590
$ws.Range("L16").Value = 2900.8
$ws.Range("M16").Value = -1420
$ws.Range("N16").Value = -3240.8

$ws.Range("H140").Value = 300104.25
$ws.Range("J140").Value = 300104.25
$ws.Range("L140").Value = 300104.25
$ws.Range("N140").Value = -310464.25
